# Adds a title block ("ResearchIt!!") and a "Team Members" list with the
# three team member names above the existing body text, per the commit
# "added title and name".

$d = $word.ActiveDocument

# Anchor on the very first paragraph of the document (the long
# "Our platform's sole purpose ..." paragraph) and insert six new
# paragraphs in front of it, in document order.
$anchor = $d.Paragraphs.First.Range

$anchor.InsertParagraphBefore()   # 1: tabs + "ResearchIt!!" title
$anchor.InsertParagraphBefore()   # 2: "Team Members:"
$anchor.InsertParagraphBefore()   # 3: "1. Raghav V Pandit - PES1UG19CS364"
$anchor.InsertParagraphBefore()   # 4: "2. Rishab Kashyap B S - PES1UG19CS385"
$anchor.InsertParagraphBefore()   # 5: "3. Royston E Tauro - PES1UG19CS396"
$anchor.InsertParagraphBefore()   # 6: blank spacer line

$colorGray = 2236962   # 0x222222

# --- Paragraph 1: leading tabs + big bold title -------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "`t`t`t`tResearchIt!!"
$p1.Range.Font.Name = "Arial"
$p1.Range.Font.Color = $colorGray
$p1.Range.Font.BoldBi = $true
$p1.Range.Font.Size = 12
$p1.Range.Font.SizeBi = 12

# Re-apply bold + larger size just over the "ResearchIt!!" text (leave the
# four leading tabs at the smaller, non-bold size).
$titleStart = $p1.Range.Start + 4
$titleRange = $d.Range($titleStart, $p1.Range.End - 1)
$titleRange.Font.Name = "Arial"
$titleRange.Font.Color = $colorGray
$titleRange.Font.Bold = $true
$titleRange.Font.BoldBi = $true
$titleRange.Font.Size = 26
$titleRange.Font.SizeBi = 26

# --- Paragraph 2: "Team Members:" ---------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = "Team Members:"
$p2.Range.Font.Name = "Arial"
$p2.Range.Font.Color = $colorGray
$p2.Range.Font.Bold = $true
$p2.Range.Font.BoldBi = $true
$p2.Range.Font.Size = 16
$p2.Range.Font.SizeBi = 16

# --- Paragraph 3: first team member -------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = "1. Raghav V Pandit - PES1UG19CS364"
$p3.Range.Font.Name = "Arial"
$p3.Range.Font.Color = $colorGray
$p3.Range.Font.Bold = $true
$p3.Range.Font.BoldBi = $true
$p3.Range.Font.Size = 16
$p3.Range.Font.SizeBi = 16

# --- Paragraph 4: second team member ------------------------------------
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Text = "2. Rishab Kashyap B S - PES1UG19CS385"
$p4.Range.Font.Name = "Arial"
$p4.Range.Font.Color = $colorGray
$p4.Range.Font.Bold = $true
$p4.Range.Font.BoldBi = $true
$p4.Range.Font.Size = 16
$p4.Range.Font.SizeBi = 16

# --- Paragraph 5: third team member --------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Text = "3. Royston E Tauro - PES1UG19CS396"
$p5.Range.Font.Name = "Arial"
$p5.Range.Font.Color = $colorGray
$p5.Range.Font.Bold = $true
$p5.Range.Font.BoldBi = $true
$p5.Range.Font.Size = 16
$p5.Range.Font.SizeBi = 16

# --- Paragraph 6: blank spacer line between the header block and body ---
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Font.Name = "Arial"
$p6.Range.Font.Color = $colorGray
$p6.Range.Font.BoldBi = $true
$p6.Range.Font.Size = 12
$p6.Range.Font.SizeBi = 12
$p6.Range.HighlightColorIndex = 8   # wdWhite

Write-Output "Inserted title/team-members block; paragraph count = $($d.Paragraphs.Count)"
